# Update "想去人数" (F column) values for both the "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 99
    4  = 1564
    5  = 599
    6  = 1087
    7  = 11301
    8  = 13
    10 = 443
    11 = 338
    12 = 1083
    13 = 777
    14 = 12299
    15 = 12942
    17 = 134
    22 = 82
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
